# 0.8 change related fixes
# - Insert a new "is_active" column into node_c (sheet "node_c") right after
#   "has_balance" (i.e. before the old column D "inflow_method"), and mark
#   is_active = "yes" for the node rows that need it.
# - Switch the active sheet/tab from "scenario" to "node_c" and update its
#   selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("node_c")

# Insert a new blank column before column D ("inflow_method"); everything
# from D onward (including the drawing's anchor) shifts one column right.
$ws.Columns("D").Insert()

# New column header + values: "is_active" flag.
$ws.Range("D2").Value = "is_active"
$ws.Range("D3").Value = "yes"
$ws.Range("D4").Value = "yes"
$ws.Range("D9").Value = "yes"

# The sheet's floating TextBox is anchored to a cell (one-cell anchor); its
# on-screen column boundary needs to be recomputed now that a column was
# inserted ahead of it, so it keeps pointing at the same visual spot.
$shp = $ws.Shapes.Item(1)
$cumWidth = 0
for ($i = 1; $i -le 34; $i++) {
  $cumWidth = $cumWidth + $ws.Cells.Item(1, $i).Width
}
$shp.Left = $cumWidth + (161925 / 12700)

# Make node_c the active/selected sheet (was "scenario" before).
$ws.Activate() | Out-Null
$ws.Range("J20").Select() | Out-Null
